$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.308.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06563"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.91"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.114"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6658"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +9.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.317.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.133.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007293"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.345"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.156"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.304"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.85"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.374"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09774"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.449"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.172"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7063"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.092"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.692"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.520"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8676"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.965"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.92"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4189"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "988.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.213"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.207"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1163"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.12%  "
